$d = $word.ActiveDocument

# 1. Fix wording: "over or above" -> "under or at least"
$d.Content.Find.Execute("over or above", $true, $false, $false, $false, $false,
                         $true, 1, $false, "under or at least", 2) | Out-Null

# 2. Fix duplicated period typo: "disease state.. We tested" -> "disease state. We tested"
$d.Content.Find.Execute("disease state.. We tested", $true, $false, $false, $false, $false,
                         $true, 1, $false, "disease state. We tested", 2) | Out-Null

# 3. Remove the embedded Figure 1 text box (image + caption) anchored in the
#    "After this analysis we found that..." paragraph.
if ($d.Shapes.Count -gt 0) {
    $d.Shapes.Item(1).Delete()
}

# 4. Simplify "genes which are enriched" -> "genes are enriched"
$d.Content.Find.Execute("genes which are enriched", $true, $false, $false, $false, $false,
                         $true, 1, $false, "genes are enriched", 2) | Out-Null

# 5. Remove the trailing sentence about Figure 1 / PCA from the "As the reviewer suspected" paragraph
$d.Content.Find.Execute("Furthermore, as shown in Figure 1 of this response, adding the age factor into our analysis dramatically improved the clustering of our data, as determined by principal component analysis, which shows that the first principal component separates older subjects from younger subjects with one exception in each case.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null
